$d = $word.ActiveDocument

# The document body is:
#   1: "TEMAT:  Plan Lekcji" (title)
#   2: (empty)
#   3: long description paragraph
#   4: (empty)   <- insert the two new paragraphs right after this one
#   5: (empty)
#   6: (empty)
#
# Insert two new empty-style paragraphs after paragraph 4, and put the
# git repo link text into the second of the two new paragraphs.

$anchor = $d.Paragraphs.Item(4)
$anchor.Range.InsertParagraphAfter()

$firstNew = $d.Paragraphs.Item(5)
$firstNew.Range.InsertParagraphAfter()

$secondNew = $d.Paragraphs.Item(6)
$secondNew.Range.Text = "git repo: https://github.com/VodkaAzFYR/DataBasesProject2023"
